# Update the 8-yue (August) sign-in sheet from the 2019-07-26..2019-08-25
# period to the 2020-08-05..2020-09-04 period, and extend the attendance
# table with the two extra working days (2020-09-03, 2020-09-04) that fall
# inside the new date range.

$d = $word.ActiveDocument

# --- 1. Replace every date token in the document (header line + all table
#        rows) with its corresponding new date. Every "old" value starts
#        with 2019- and every "new" value starts with 2020-, so a simple
#        sequence of whole-document Find/Replace-All calls is unambiguous
#        and order-independent.
$dateMap = @(
    @("2019-07-26", "2020-08-05"),
    @("2019-08-25", "2020-09-04"),
    @("2019-07-29", "2020-08-06"),
    @("2019-07-30", "2020-08-07"),
    @("2019-07-31", "2020-08-10"),
    @("2019-08-01", "2020-08-11"),
    @("2019-08-02", "2020-08-12"),
    @("2019-08-05", "2020-08-13"),
    @("2019-08-06", "2020-08-14"),
    @("2019-08-07", "2020-08-17"),
    @("2019-08-08", "2020-08-18"),
    @("2019-08-09", "2020-08-19"),
    @("2019-08-12", "2020-08-20"),
    @("2019-08-13", "2020-08-21"),
    @("2019-08-14", "2020-08-24"),
    @("2019-08-15", "2020-08-25"),
    @("2019-08-16", "2020-08-26"),
    @("2019-08-19", "2020-08-27"),
    @("2019-08-20", "2020-08-28"),
    @("2019-08-21", "2020-08-31"),
    @("2019-08-22", "2020-09-01"),
    @("2019-08-23", "2020-09-02")
)

foreach ($pair in $dateMap) {
    $oldDate = $pair[0]
    $newDate = $pair[1]
    $d.Content.Find.Execute($oldDate, $true, $false, $false, $false, $false,
                             $true, 1, $false, $newDate, 2)
}

# --- 2. Append two new attendance rows (index 22 / 2020-09-03 and index 23
#        / 2020-09-04) to the bottom of the sign-in table, matching the
#        layout of the existing rows (index + date filled in, the
#        sign-in/sign-out/remarks cells left blank).
$t = $d.Tables.Item(1)

$newRows = @(
    @("22", "2020-09-03"),
    @("23", "2020-09-04")
)

foreach ($row in $newRows) {
    $t.Rows.Add() | Out-Null
    $rowIndex = $t.Rows.Count
    $t.Cell($rowIndex, 1).Range.Text = $row[0]
    $t.Cell($rowIndex, 2).Range.Text = $row[1]
}

Write-Output "done"
